$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (49) down through the new rows (50-74)
$ws.Range("A49:H49").Copy()
$ws.Range("A50:H74").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 50
$ws.Cells.Item(50, 1).Value = 'Wk28'
$ws.Cells.Item(50, 2).Value = 'SAP'
$ws.Cells.Item(50, 3).Value = 'XS-MIT-0099'
$ws.Cells.Item(50, 4).Value = 'SY3120-5LZDM5-F2 SOLENOID VALVE'
$ws.Cells.Item(50, 5).Value = '1pcs'
$ws.Cells.Item(50, 6).Value = 'Fishes'
$ws.Cells.Item(50, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 7).Date
$ws.Cells.Item(50, 8).Value = 74.87

# Row 51
$ws.Cells.Item(51, 1).Value = 'Wk28'
$ws.Cells.Item(51, 2).Value = 'SAP'
$ws.Cells.Item(51, 3).Value = 11155497
$ws.Cells.Item(51, 4).Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Cells.Item(51, 5).Value = '12pcs'
$ws.Cells.Item(51, 6).Value = 'Fishes'
$ws.Cells.Item(51, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 7).Date
$ws.Cells.Item(51, 8).Value = 421.2

# Row 52
$ws.Cells.Item(52, 1).Value = 'Wk28'
$ws.Cells.Item(52, 2).Value = 'SAP'
$ws.Cells.Item(52, 3).Value = 11155497
$ws.Cells.Item(52, 4).Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Cells.Item(52, 5).Value = '12pcs'
$ws.Cells.Item(52, 6).Value = 'Fishes'
$ws.Cells.Item(52, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(52, 8).Value = 421.2

# Row 53
$ws.Cells.Item(53, 1).Value = 'Wk28'
$ws.Cells.Item(53, 2).Value = 'SAP'
$ws.Cells.Item(53, 3).Value = 'XS-PTS-1030'
$ws.Cells.Item(53, 4).Value = 'TW.50.1A.00.02.016.00 BRASS SLIDE'
$ws.Cells.Item(53, 5).Value = '10pcs'
$ws.Cells.Item(53, 6).Value = 'Fishes'
$ws.Cells.Item(53, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(53, 8).Value = 13.53

# Row 54
$ws.Cells.Item(54, 1).Value = 'Wk28'
$ws.Cells.Item(54, 2).Value = 'SAP'
$ws.Cells.Item(54, 3).Value = 'XS-PTS-0860'
$ws.Cells.Item(54, 4).Value = 'X1439 70192.312 DARK LEAKAGE RUBBER TIP'
$ws.Cells.Item(54, 5).Value = '20pcs'
$ws.Cells.Item(54, 6).Value = 'Fishes'
$ws.Cells.Item(54, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(54, 8).Value = 257.2

# Row 55
$ws.Cells.Item(55, 1).Value = 'Wk28'
$ws.Cells.Item(55, 2).Value = 'SAP'
$ws.Cells.Item(55, 3).Value = 'XS-PTS-0764'
$ws.Cells.Item(55, 4).Value = '70192.384 RUBBER TIP'
$ws.Cells.Item(55, 5).Value = '10pcs'
$ws.Cells.Item(55, 6).Value = 'Fishes'
$ws.Cells.Item(55, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(55, 8).Value = 214.3

# Row 56
$ws.Cells.Item(56, 1).Value = 'Wk28'
$ws.Cells.Item(56, 2).Value = 'SAP'
$ws.Cells.Item(56, 3).Value = 'XS-PTS-0790'
$ws.Cells.Item(56, 4).Value = '70900.108 Test Socket Pogo Pin'
$ws.Cells.Item(56, 5).Value = '200pcs'
$ws.Cells.Item(56, 6).Value = 'Fishes'
$ws.Cells.Item(56, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(56, 8).Value = 896

# Row 57
$ws.Cells.Item(57, 1).Value = 'Wk28'
$ws.Cells.Item(57, 2).Value = 'SAP'
$ws.Cells.Item(57, 3).Value = 'XS-PTS-0933'
$ws.Cells.Item(57, 4).Value = 'X1629 70900.137 SEMICONDUCTOR PROBE PIN'
$ws.Cells.Item(57, 5).Value = '200pcs'
$ws.Cells.Item(57, 6).Value = 'Fishes'
$ws.Cells.Item(57, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(57, 8).Value = 2482

# Row 58
$ws.Cells.Item(58, 1).Value = 'Wk28'
$ws.Cells.Item(58, 2).Value = 'SAP'
$ws.Cells.Item(58, 3).Value = 'XS-PTS-0876'
$ws.Cells.Item(58, 4).Value = 'HX 2067 HPN CRC CO CONTACT CLEANER'
$ws.Cells.Item(58, 5).Value = '4pcs'
$ws.Cells.Item(58, 6).Value = 'Fishes'
$ws.Cells.Item(58, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(58, 8).Value = 91.88

# Row 59
$ws.Cells.Item(59, 1).Value = 'Wk28'
$ws.Cells.Item(59, 2).Value = 'SAP'
$ws.Cells.Item(59, 3).Value = 'XS-SPM-0061'
$ws.Cells.Item(59, 4).Value = '14210873.140 SPAREP A218 - PICKUP'
$ws.Cells.Item(59, 5).Value = '3pcs'
$ws.Cells.Item(59, 6).Value = 'Lisa'
$ws.Cells.Item(59, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 8).Date
$ws.Cells.Item(59, 8).Value = 621.36

# Row 60
$ws.Cells.Item(60, 1).Value = 'Wk28'
$ws.Cells.Item(60, 2).Value = 'SAP'
$ws.Cells.Item(60, 3).Value = 'XS-SPE-0090'
$ws.Cells.Item(60, 4).Value = 'GP-762A35A-02XB Hyperspace Semicon Probe'
$ws.Cells.Item(60, 5).Value = '200pcs'
$ws.Cells.Item(60, 6).Value = 'Lisa'
$ws.Cells.Item(60, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 9).Date
$ws.Cells.Item(60, 8).Value = 312

# Row 61
$ws.Cells.Item(61, 1).Value = 'Wk28'
$ws.Cells.Item(61, 2).Value = 'SAP'
$ws.Cells.Item(61, 3).Value = 'XS-PTS-0356'
$ws.Cells.Item(61, 4).Value = 'GP-570D84A-03 Hyperspace Semicon Probes'
$ws.Cells.Item(61, 5).Value = '200pcs'
$ws.Cells.Item(61, 6).Value = 'Sihl'
$ws.Cells.Item(61, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 9).Date
$ws.Cells.Item(61, 8).Value = 420

# Row 62
$ws.Cells.Item(62, 1).Value = 'Wk28'
$ws.Cells.Item(62, 2).Value = 'SAP'
$ws.Cells.Item(62, 3).Value = 'XS-PTS-0864'
$ws.Cells.Item(62, 4).Value = 'Model : 6K-76235-H03X-NST Hyperspace'
$ws.Cells.Item(62, 5).Value = '1pcs'
$ws.Cells.Item(62, 6).Value = 'Sihl'
$ws.Cells.Item(62, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 9).Date
$ws.Cells.Item(62, 8).Value = 325.03

# Row 63
$ws.Cells.Item(63, 1).Value = 'Wk28'
$ws.Cells.Item(63, 2).Value = 'SAP'
$ws.Cells.Item(63, 3).Value = 'XS-PTS-0868'
$ws.Cells.Item(63, 4).Value = '6K-76235-H03X-ALN Hyperspace Semicon'
$ws.Cells.Item(63, 5).Value = '1pcs'
$ws.Cells.Item(63, 6).Value = 'Sihl'
$ws.Cells.Item(63, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 9).Date
$ws.Cells.Item(63, 8).Value = 210.07

# Row 64
$ws.Cells.Item(64, 1).Value = 'Wk28'
$ws.Cells.Item(64, 2).Value = 'SAP'
$ws.Cells.Item(64, 3).Value = 11155497
$ws.Cells.Item(64, 4).Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Cells.Item(64, 5).Value = '12pcs'
$ws.Cells.Item(64, 6).Value = 'Fishes'
$ws.Cells.Item(64, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 9).Date
$ws.Cells.Item(64, 8).Value = 421.2

# Row 65
$ws.Cells.Item(65, 1).Value = 'Wk28'
$ws.Cells.Item(65, 2).Value = 'SAP'
$ws.Cells.Item(65, 3).NumberFormat = "@"
$ws.Cells.Item(65, 3).Value = '11151237'
$ws.Cells.Item(65, 4).Value = 'PTS-1060 70192.692 X2637 RUBBER TIP'
$ws.Cells.Item(65, 5).Value = '60pcs'
$ws.Cells.Item(65, 6).Value = 'Fishes'
$ws.Cells.Item(65, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(65, 8).Value = 406.2

# Row 66
$ws.Cells.Item(66, 1).Value = 'Wk28'
$ws.Cells.Item(66, 2).Value = 'SAP'
$ws.Cells.Item(66, 3).NumberFormat = "@"
$ws.Cells.Item(66, 3).Value = '11155143'
$ws.Cells.Item(66, 4).Value = 'PTS-1136 300-001519-015 Semiconductor'
$ws.Cells.Item(66, 5).Value = '5200pcs'
$ws.Cells.Item(66, 6).Value = 'Fishes'
$ws.Cells.Item(66, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(66, 8).Value = 9672

# Row 67
$ws.Cells.Item(67, 1).Value = 'Wk28'
$ws.Cells.Item(67, 2).Value = 'SAP'
$ws.Cells.Item(67, 3).Value = 'XS-PTS-0469'
$ws.Cells.Item(67, 4).Value = '70508.201 10# BULB MFG: PHILIPS'
$ws.Cells.Item(67, 5).Value = '2pcs'
$ws.Cells.Item(67, 6).Value = 'Fishes'
$ws.Cells.Item(67, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(67, 8).Value = 263.76

# Row 68
$ws.Cells.Item(68, 1).Value = 'Wk28'
$ws.Cells.Item(68, 2).Value = 'SAP'
$ws.Cells.Item(68, 3).Value = 'XS-PTS-0465'
$ws.Cells.Item(68, 4).Value = '70508.190 00707-SU Bulb#6'
$ws.Cells.Item(68, 5).Value = '3pcs'
$ws.Cells.Item(68, 6).Value = 'Fishes'
$ws.Cells.Item(68, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(68, 8).Value = 319.26

# Row 69
$ws.Cells.Item(69, 1).Value = 'Wk28'
$ws.Cells.Item(69, 2).Value = 'SAP'
$ws.Cells.Item(69, 3).Value = 11155497
$ws.Cells.Item(69, 4).Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Cells.Item(69, 5).Value = '12pcs'
$ws.Cells.Item(69, 6).Value = 'Fishes'
$ws.Cells.Item(69, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(69, 8).Value = 421.2

# Row 70
$ws.Cells.Item(70, 1).Value = 'Wk28'
$ws.Cells.Item(70, 2).Value = 'SAP'
$ws.Cells.Item(70, 3).Value = 'XS-PTS-0356'
$ws.Cells.Item(70, 4).Value = 'GP-570D84A-03 Hyperspace Semicon Probes'
$ws.Cells.Item(70, 5).Value = '200pcs'
$ws.Cells.Item(70, 6).Value = 'Sihl'
$ws.Cells.Item(70, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(70, 8).Value = 420

# Row 71
$ws.Cells.Item(71, 1).Value = 'Wk28'
$ws.Cells.Item(71, 2).Value = 'SAP'
$ws.Cells.Item(71, 3).Value = 'XS-PTS-0864'
$ws.Cells.Item(71, 4).Value = 'Model : 6K-76235-H03X-NST Hyperspace'
$ws.Cells.Item(71, 5).Value = '1pcs'
$ws.Cells.Item(71, 6).Value = 'Sihl'
$ws.Cells.Item(71, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(71, 8).Value = 325.03

# Row 72
$ws.Cells.Item(72, 1).Value = 'Wk28'
$ws.Cells.Item(72, 2).Value = 'SAP'
$ws.Cells.Item(72, 3).Value = 'XS-PTS-0868'
$ws.Cells.Item(72, 4).Value = '6K-76235-H03X-ALN Hyperspace Semicon'
$ws.Cells.Item(72, 5).Value = '1pcs'
$ws.Cells.Item(72, 6).Value = 'Sihl'
$ws.Cells.Item(72, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(72, 8).Value = 210.07

# Row 73
$ws.Cells.Item(73, 1).Value = 'Wk28'
$ws.Cells.Item(73, 2).Value = 'SAP'
$ws.Cells.Item(73, 3).Value = 11155497
$ws.Cells.Item(73, 4).Value = 'PTS-1138 6K-9090-HS01  Vacuum Adaptor'
$ws.Cells.Item(73, 5).Value = '8pcs'
$ws.Cells.Item(73, 6).Value = 'Fishes'
$ws.Cells.Item(73, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 10).Date
$ws.Cells.Item(73, 8).Value = 280.8

# Row 74
$ws.Cells.Item(74, 1).Value = 'Wk28'
$ws.Cells.Item(74, 2).Value = 'SAP'
$ws.Cells.Item(74, 3).Value = 11155797
$ws.Cells.Item(74, 4).Value = 'PTS-1140 6K-9090-HS02 HPS LF Aligner'
$ws.Cells.Item(74, 5).Value = '18pcs'
$ws.Cells.Item(74, 6).Value = 'Fishes'
$ws.Cells.Item(74, 7).Value = (Get-Date -Year 2025 -Month 7 -Day 11).Date
$ws.Cells.Item(74, 8).Value = 362.88

# Restore default (General/style 1) formatting on cells forced to text so the
# pasted number style doesn't linger (matches style of other text-number cells like C5)
$ws.Range("C5").Copy()
$ws.Range("C65").PasteSpecial(-4122)
$ws.Range("C66").PasteSpecial(-4122)
$excel.CutCopyMode = $false
